# Correctly writes the winners of the first round of the first bracket.
#
# - On the "Bracket" sheet, the D-column "winner" slot for each of the 8
#   first-round matchups in the first (left-hand) bracket is filled in with
#   the chosen winner's name, replacing the placeholder "anN" labels that
#   used to sit in D/N for that row; the bracket-position number in column A
#   and the mirrored label in column N for that row are cleared out too.
# - A brand-new "Bracket1" sheet is added (right after "Bracket") that holds
#   just those same eight round-1 picks in its D column, recording the
#   winners of the first round of the first bracket on their own.
# - "Bracket" ends up the selected/active sheet.

$wb = $excel.ActiveWorkbook
$bracketSheet = $wb.Worksheets.Item("Bracket")

function Set-Winner($Sheet, $Row, $Winner, $ClearSides) {
    if ($ClearSides) {
        $Sheet.Range("A$Row").ClearContents()
        $Sheet.Range("N$Row").ClearContents()
    }
    $Sheet.Range("D$Row").Value = $Winner
}

# Round-1 winners, first bracket (rows 2,6,10,...,30), as currently recorded
# on the "Bracket" sheet.
Set-Winner $bracketSheet 2  "Okapi" $true
Set-Winner $bracketSheet 6  "Striped Polecat" $true
Set-Winner $bracketSheet 10 "Side-striped jackal" $true
Set-Winner $bracketSheet 14 "Striped dolphin" $true
Set-Winner $bracketSheet 18 "Wildcat" $true
Set-Winner $bracketSheet 22 "Striped hyena" $true
Set-Winner $bracketSheet 26 "Striped Rabbit" $true
Set-Winner $bracketSheet 30 "Kudu" $true

# New "Bracket1" sheet, placed right after "Bracket", holding the same
# round-1 picks for the first bracket.
$bracket1Sheet = $wb.Worksheets.Add($null, $bracketSheet)
$bracket1Sheet.Name = "Bracket1"

Set-Winner $bracket1Sheet 2  "Okapi" $false
Set-Winner $bracket1Sheet 6  "Striped Polecat" $false
Set-Winner $bracket1Sheet 10 "Side-striped jackal" $false
Set-Winner $bracket1Sheet 14 "Striped dolphin" $false
Set-Winner $bracket1Sheet 18 "Wildcat" $false
Set-Winner $bracket1Sheet 22 "Striped hyena" $false
Set-Winner $bracket1Sheet 26 "Numbat" $false
Set-Winner $bracket1Sheet 30 "Kudu" $false

# Leave "Bracket" as the selected tab.
$bracketSheet.Activate()
